$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) "Förändrad" (changed/updated) date in column C moved forward two days
#    (45184 -> 45186) for every data row (2 through 130).
$ws.Range("C2:C130").Value = 45186

# 2) For the rows that carry species-link columns (S, T, V, W, X, Y), the
#    HYPERLINK() formulas gain a friendly-name second argument equal to the
#    "Beteckning" (case id) shown in column A for that row.
$linkCols = @(
    @{ Col = "S"; Path = "artfynd";         Ext = "xlsx" },
    @{ Col = "T"; Path = "kartor";          Ext = "png"  },
    @{ Col = "V"; Path = "klagomål";        Ext = "docx" },
    @{ Col = "W"; Path = "klagomålsmail";   Ext = "docx" },
    @{ Col = "X"; Path = "tillsyn";         Ext = "docx" },
    @{ Col = "Y"; Path = "tillsynsmail";    Ext = "docx" }
)

for ($row = 2; $row -le 47; $row++) {
    $beteckning = $ws.Cells.Item($row, 1).Text

    foreach ($entry in $linkCols) {
        $col = $entry.Col
        $path = $entry.Path
        $ext = $entry.Ext
        $cellRef = "$col$row"
        $url = "https://klasma.github.io/Logging_BORGHOLM/$path/$beteckning.$ext"
        $formula = '=HYPERLINK("' + $url + '", "' + $beteckning + '")'
        $ws.Range($cellRef).Formula = $formula
    }
}
